$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# -----------------------------------------------------------------
# 1) Duplicate the "Week 10" block (rows 69-76) down to rows 77-84
#    while it is still in its original (pre-edit) state, because the
#    new "Week 12" block re-uses that original content / formatting.
# -----------------------------------------------------------------
$srcBlock = $ws.Range("A69:G76")
$dstBlock = $ws.Range("A77:G84")
$srcBlock.Copy($dstBlock)
$excel.CutCopyMode = $false

# Rows 80-83 should not contain any Actual-Hours (column E) cell at all
# (mirrors the original un-edited rows 72-75 that fed them).
$ws.Cells.Item(80,5).Clear()
$ws.Cells.Item(81,5).Clear()
$ws.Cells.Item(82,5).Clear()
$ws.Cells.Item(83,5).Clear()

# Row 79 gets a taller height (wrapped "Team meeting" text)
$ws.Rows.Item(79).RowHeight = 20.4

# B80/C80 need to switch from the plain style to the shaded style
# already used on B71/C71 (which also already hold "All"/"All").
$ws.Cells.Item(71,2).Copy()
$ws.Cells.Item(80,2).PasteSpecial(-4122)
$ws.Cells.Item(71,3).Copy()
$ws.Cells.Item(80,3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# 2) Fill in the new activities for rows 80-83 (new "Week 12" rows)
# -----------------------------------------------------------------
$ws.Cells.Item(80,1).Value2 = "Client Office visit"
$ws.Cells.Item(81,1).Value2 = "Status Tracker Page update"
$ws.Cells.Item(83,1).Value2 = "Start working on R2"

$ws.Cells.Item(80,3).Value2 = "All"
$ws.Cells.Item(80,4).Value2 = 4

$ws.Cells.Item(81,3).Value2 = "Agrim"
$ws.Cells.Item(81,4).Value2 = 2

$ws.Cells.Item(82,1).Value2 = "Live testing"
$ws.Cells.Item(82,4).Value2 = 1

$ws.Cells.Item(83,4).Value2 = 2

$ws.Cells.Item(80,7).Value2 = "Visit the client office and get insights of project"
$ws.Cells.Item(81,7).Value2 = "Add few options in status tracker page"
$ws.Cells.Item(82,7).Value2 = "Check all the purchase functionality"
$ws.Cells.Item(83,7).Value2 = "Plan and divide each work for R2"

# F82/F83 keep the "Planned" text but switch to the plain (unshaded) style
$ws.Cells.Item(72,1).Copy()
$ws.Cells.Item(82,6).PasteSpecial(-4122)
$ws.Cells.Item(83,6).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "Week 12" banner row
$ws.Cells.Item(84,1).Value2 = "Week 12  (April 01-April 07  )"

# -----------------------------------------------------------------
# 3) Update the original "Week 10" rows (69-75): fill Actual Hours
#    (column E) and flip the Status (column F) to Done / Ongoing.
# -----------------------------------------------------------------
$ws.Cells.Item(69,5).Value2 = 0.5
$ws.Cells.Item(70,5).Value2 = 0.5
$ws.Cells.Item(71,5).Value2 = 2
$ws.Cells.Item(72,5).Value2 = 3
$ws.Cells.Item(73,5).Value2 = 2
$ws.Cells.Item(74,5).Value2 = 1

# F69:F74 -> "Done" status, re-using the green "Done" styling already
# present on F62 (style index 32).
$ws.Cells.Item(62,6).Copy()
$ws.Cells.Item(69,6).PasteSpecial(-4122)
$ws.Cells.Item(70,6).PasteSpecial(-4122)
$ws.Cells.Item(71,6).PasteSpecial(-4122)
$ws.Cells.Item(72,6).PasteSpecial(-4122)
$ws.Cells.Item(73,6).PasteSpecial(-4122)
$ws.Cells.Item(74,6).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(69,6).Value2 = "Done"
$ws.Cells.Item(70,6).Value2 = "Done"
$ws.Cells.Item(71,6).Value2 = "Done"
$ws.Cells.Item(72,6).Value2 = "Done"
$ws.Cells.Item(73,6).Value2 = "Done"
$ws.Cells.Item(74,6).Value2 = "Done"

# F75 -> "Ongoing" status, re-using the red "Ongoing" styling already
# present on F64 (style index 31).
$ws.Cells.Item(64,6).Copy()
$ws.Cells.Item(75,6).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(75,6).Value2 = "Ongoing"

# -----------------------------------------------------------------
# 4) Move the active selection down to the new first blank row below
#    the freshly added "Week 12" block, and scroll the view down.
# -----------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("B87").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 66
